$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 21500
$ws.Range("H23").Value = 21500
$ws.Range("H33").Value = 100.7
$ws.Range("I33").Value = 100.7
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 100.7
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = 128.3
$ws.Range("H74").Value = 12875
$ws.Range("I74").Value = 4791.6665
$ws.Range("K74").Value = 4791.6665
$ws.Range("M74").Value = -3855.6665
$ws.Range("H77").Value = 12875
$ws.Range("I77").Value = 4791.6665
$ws.Range("K77").Value = 23958.3325
$ws.Range("M77").Value = -19278.3325
$ws.Range("H87").Value = 59538
$ws.Range("J87").Value = 59538
$ws.Range("L87").Value = 59538
$ws.Range("N87").Value = -62034
$ws.Range("H90").Value = 59538
$ws.Range("J90").Value = 59538
$ws.Range("L90").Value = 178614
$ws.Range("N90").Value = -191094
$ws.Range("H96").Value = 558.9231
$ws.Range("I96").Value = 359.1
$ws.Range("K96").Value = 1077.3
$ws.Range("M96").Value = 295.6999999999998
$ws.Range("H112").Value = 2700
$ws.Range("J112").Value = 2875
$ws.Range("L112").Value = 8625
$ws.Range("N112").Value = -10841
$ws.Range("H113").Value = 1399
$ws.Range("I113").Value = 1399
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1399
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 1855
$ws.Range("H121").Value = 820.6
$ws.Range("J121").Value = 820.6
$ws.Range("L121").Value = 2461.8
$ws.Range("N121").Value = -5955.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 111
$ws.Range("I5").Value = 111
$ws.Range("K5").Value = 111
$ws.Range("M5").Value = 1
$ws.Range("H32").Value = 5240.8
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").ClearContents()
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = 0
$ws.Range("H74").Value = 4725.476
$ws.Range("I74").Value = 3337.8572
$ws.Range("K74").Value = 3337.8572
$ws.Range("M74").Value = -2463.8572
$ws.Range("H77").Value = 4725.476
$ws.Range("I77").Value = 3337.8572
$ws.Range("K77").Value = 16689.286
$ws.Range("M77").Value = -12321.286
$ws.Range("H102").Value = 7624.75
$ws.Range("I102").Value = 5499.5
$ws.Range("J102").Value = 9750
$ws.Range("K102").Value = 5499.5
$ws.Range("L102").Value = 9750
$ws.Range("M102").Value = -3877.5
$ws.Range("N102").Value = -12994
$ws.Range("H132").Value = 5149.6665
$ws.Range("I132").Value = 1966.3334
$ws.Range("J132").Value = 8333
$ws.Range("K132").Value = 5899.0002
$ws.Range("L132").Value = 24999
$ws.Range("M132").Value = -3369.0002
$ws.Range("N132").Value = -30059
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = 0
$ws.Range("H141").Value = 150000
$ws.Range("J141").Value = 150000
$ws.Range("L141").Value = 150000
$ws.Range("N141").Value = -160360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 111
$ws.Range("I4").Value = 111
$ws.Range("K4").Value = 111
$ws.Range("M4").Value = 4
$ws.Range("H15").Value = 2500
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 2500
$ws.Range("K15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("M15").Value = 2500
$ws.Range("N15").Value = -2954
$ws.Range("H35").Value = 62445
$ws.Range("J35").Value = 62445
$ws.Range("L35").Value = 62445
$ws.Range("N35").Value = -63065
$ws.Range("H107").Value = 4182.143
$ws.Range("I107").Value = 1393.2727
$ws.Range("K107").Value = 1393.2727
$ws.Range("M107").Value = 526.7273
$ws.Range("H134").Value = 3525.7778
$ws.Range("I134").Value = 3059.4285
$ws.Range("K134").Value = 9178.2855
$ws.Range("M134").Value = -6643.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 71.63636
$ws.Range("I7").Value = 49.375
$ws.Range("J7").Value = 131
$ws.Range("K7").Value = 49.375
$ws.Range("L7").Value = 131
$ws.Range("M7").Value = 63.625
$ws.Range("N7").Value = -357
$ws.Range("H16").Value = 1435.9
$ws.Range("I16").Value = 1401.1111
$ws.Range("K16").Value = 1401.1111
$ws.Range("M16").Value = -1114.1111
$ws.Range("H22").Value = 1293.2
$ws.Range("I22").Value = 366.5
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 366.5
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -16.5
$ws.Range("N22").Value = -5700
$ws.Range("H33").Value = 1307.2858
$ws.Range("I33").Value = 1191.8334
$ws.Range("K33").Value = 1191.8334
$ws.Range("M33").Value = -812.8334
$ws.Range("H86").Value = 1913.5
$ws.Range("I86").Value = 1896.2
$ws.Range("K86").Value = 1896.2
$ws.Range("M86").Value = -773.2
$ws.Range("H89").Value = 1913.5
$ws.Range("I89").Value = 1896.2
$ws.Range("K89").Value = 9481
$ws.Range("M89").Value = -3865
$ws.Range("H113").Value = 1435.9
$ws.Range("I113").Value = 1401.1111
$ws.Range("K113").Value = 1401.1111
$ws.Range("M113").Value = 768.8888999999999
$ws.Range("H115").Value = 75332.336
$ws.Range("J115").Value = 75332.336
$ws.Range("L115").Value = 75332.336
$ws.Range("N115").Value = -77682.336
$ws.Range("H122").Value = 1423.5714
$ws.Range("I122").Value = 1274.2222
$ws.Range("K122").Value = 3822.6666
$ws.Range("M122").Value = -1372.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 41.666668
$ws.Range("I7").Value = 37
$ws.Range("J7").Value = 51
$ws.Range("K7").Value = 111
$ws.Range("L7").Value = 153
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = -377
$ws.Range("H8").Value = 387.8
$ws.Range("I8").Value = 387.8
$ws.Range("K8").Value = 1163.4
$ws.Range("M8").Value = -1024.4
$ws.Range("H44").Value = 281.29413
$ws.Range("I44").Value = 101.64286
$ws.Range("K44").Value = 304.92858
$ws.Range("M44").Value = 93.07141999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2033
$ws.Range("I80").Value = 1874.75
$ws.Range("J80").Value = 2159.6
$ws.Range("K80").Value = 1874.75
$ws.Range("L80").Value = 2159.6
$ws.Range("M80").Value = -876.75
$ws.Range("N80").Value = -4155.6
$ws.Range("H83").Value = 2033
$ws.Range("I83").Value = 1874.75
$ws.Range("J83").Value = 2159.6
$ws.Range("K83").Value = 9373.75
$ws.Range("L83").Value = 10798
$ws.Range("M83").Value = -4381.75
$ws.Range("N83").Value = -20782
$ws.Range("H113").Value = 8127
$ws.Range("J113").Value = 9078
$ws.Range("L113").Value = 9078
$ws.Range("N113").Value = -13418

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1778.2858
$ws.Range("I22").Value = 1349.8
$ws.Range("J22").Value = 2849.5
$ws.Range("K22").Value = 1349.8
$ws.Range("L22").Value = 2849.5
$ws.Range("M22").Value = -1054.8
$ws.Range("N22").Value = -3439.5
$ws.Range("H27").Value = 1778.2858
$ws.Range("I27").Value = 1349.8
$ws.Range("J27").Value = 2849.5
$ws.Range("K27").Value = 1349.8
$ws.Range("L27").Value = 2849.5
$ws.Range("M27").Value = -1242.8
$ws.Range("N27").Value = -3063.5
$ws.Range("H40").Value = 7098.7856
$ws.Range("I40").Value = 6264.2
$ws.Range("K40").Value = 6264.2
$ws.Range("M40").Value = -6128.2
$ws.Range("H46").Value = 4272.778
$ws.Range("I46").Value = 2763.1428
$ws.Range("J46").Value = 5233.4546
$ws.Range("K46").Value = 2763.1428
$ws.Range("L46").Value = 5233.4546
$ws.Range("M46").Value = -2575.1428
$ws.Range("N46").Value = -5609.4546
$ws.Range("H55").Value = 1081.4546
$ws.Range("I55").Value = 1588.8334
$ws.Range("J55").Value = 472.6
$ws.Range("K55").Value = 1588.8334
$ws.Range("L55").Value = 472.6
$ws.Range("M55").Value = -1415.8334
$ws.Range("N55").Value = -818.6
$ws.Range("H61").Value = 7197.9
$ws.Range("I61").Value = 5997.25
$ws.Range("J61").Value = 7998.3335
$ws.Range("K61").Value = 5997.25
$ws.Range("L61").Value = 7998.3335
$ws.Range("M61").Value = -5795.25
$ws.Range("N61").Value = -8402.333500000001
$ws.Range("H68").Value = 8340.200000000001
$ws.Range("J68").Value = 8340.200000000001
$ws.Range("L68").Value = 8340.200000000001
$ws.Range("N68").Value = -9838.200000000001
$ws.Range("H71").Value = 8340.200000000001
$ws.Range("J71").Value = 8340.200000000001
$ws.Range("L71").Value = 41701
$ws.Range("N71").Value = -49189
$ws.Range("H113").Value = 7197.9
$ws.Range("I113").Value = 5997.25
$ws.Range("J113").Value = 7998.3335
$ws.Range("K113").Value = 5997.25
$ws.Range("L113").Value = 7998.3335
$ws.Range("M113").Value = -3827.25
$ws.Range("N113").Value = -12338.3335
$ws.Range("H136").Value = 3911.3333
$ws.Range("I136").Value = 3725.182
$ws.Range("K136").Value = 11175.546
$ws.Range("M136").Value = -8625.545999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H132").Value = 2674.75
$ws.Range("I132").Value = 2574.6924
$ws.Range("K132").Value = 7724.0772
$ws.Range("M132").Value = -5194.0772
